# order persistence and customer ui
# Append the first persisted order row (T-100) under the existing header
# row on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "T-100"
$ws.Range("B2").Value = "NTU"
$ws.Range("C2").Value = "Chicken tenders"

# "customisations" is blank for this order, but still needs a literal
# (shared-string) empty value rather than a cleared cell, and "status"
# holds the literal text "true" rather than the boolean TRUE - both are
# achieved with a leading apostrophe to force text entry, then the
# style is reset to Normal so no stray "quote prefix" cell format lingers.
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'true"
$ws.Range("D2:E2").Style = "Normal"

$ws.Range("F2").Value = "NEW"
$ws.Range("G2").Value = "Cash"
